$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 54381
$ws.Range("J3").Value = 54381
$ws.Range("L3").Value = 54381
$ws.Range("N3").Value = -54609
$ws.Range("H8").Value = 277.93332
$ws.Range("I8").Value = 62
$ws.Range("J8").Value = 299
$ws.Range("K8").Value = 186
$ws.Range("L8").Value = 897
$ws.Range("M8").Value = -47
$ws.Range("N8").Value = -1175
$ws.Range("H9").Value = 781.3333
$ws.Range("I9").Value = 1685
$ws.Range("J9").Value = 329.5
$ws.Range("K9").Value = 1685
$ws.Range("L9").Value = 329.5
$ws.Range("M9").Value = -1516
$ws.Range("N9").Value = -667.5
$ws.Range("H32").Value = 1997.25
$ws.Range("I32").Value = 1993
$ws.Range("K32").Value = 1993
$ws.Range("M32").Value = -1667
$ws.Range("H76").Value = 31256086
$ws.Range("I76").Value = 6075
$ws.Range("K76").Value = 6075
$ws.Range("M76").Value = -5760
$ws.Range("H79").Value = 31256086
$ws.Range("I79").Value = 6075
$ws.Range("K79").Value = 6075
$ws.Range("M79").Value = -4983
$ws.Range("H96").Value = 1185.6
$ws.Range("I96").Value = 983.25
$ws.Range("K96").Value = 2949.75
$ws.Range("M96").Value = -1576.75
$ws.Range("H102").Value = 54381
$ws.Range("J102").Value = 54381
$ws.Range("L102").Value = 54381
$ws.Range("N102").Value = -60871
$ws.Range("H103").Value = 799.25
$ws.Range("I103").Value = 572.3
$ws.Range("J103").Value = 961.3570999999999
$ws.Range("K103").Value = 1716.9
$ws.Range("L103").Value = 2884.0713
$ws.Range("M103").Value = -1130.9
$ws.Range("N103").Value = -4056.0713
$ws.Range("H119").Value = 1
$ws.Range("J119").Value = 1
$ws.Range("L119").Value = 3
$ws.Range("N119").Value = -9679
$ws.Range("H132").Value = 858
$ws.Range("I132").Value = 852.7778
$ws.Range("K132").Value = 2558.3334
$ws.Range("M132").Value = -28.33339999999998
$ws.Range("H138").Value = 2385962.5
$ws.Range("I138").Value = 2579.6667
$ws.Range("J138").Value = 3710064.2
$ws.Range("K138").Value = 7739.000100000001
$ws.Range("L138").Value = 11130192.6
$ws.Range("M138").Value = -2599.000100000001
$ws.Range("N138").Value = -11140472.6

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3642600
$ws.Range("I32").Value = 3850634.8
$ws.Range("K32").Value = 3850634.8
$ws.Range("M32").Value = -3850347.8
$ws.Range("H57").Value = 4999.1665
$ws.Range("I57").Value = 4999.1665
$ws.Range("K57").Value = 4999.1665
$ws.Range("M57").Value = -4515.1665
$ws.Range("H122").Value = 4599.12
$ws.Range("I122").Value = 2620.6924
$ws.Range("J122").Value = 6742.4165
$ws.Range("K122").Value = 7862.0772
$ws.Range("L122").Value = 20227.2495
$ws.Range("M122").Value = -5412.0772
$ws.Range("N122").Value = -25127.2495

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 7069663
$ws.Range("I22").Value = 9259534
$ws.Range("J22").Value = 500049.5
$ws.Range("K22").Value = 9259534
$ws.Range("L22").Value = 500049.5
$ws.Range("M22").Value = -9259361
$ws.Range("N22").Value = -500395.5
$ws.Range("H68").Value = 70000
$ws.Range("J68").Value = 70000
$ws.Range("L68").Value = 70000
$ws.Range("N68").Value = -71622
$ws.Range("H71").Value = 70000
$ws.Range("J71").Value = 70000
$ws.Range("L71").Value = 210000
$ws.Range("N71").Value = -218112
$ws.Range("H128").Value = 3862.1428
$ws.Range("I128").Value = 3862.1428
$ws.Range("K128").Value = 11586.4284
$ws.Range("M128").Value = -9096.428400000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H9").Value = 150000
$ws.Range("J9").Value = 150000
$ws.Range("L9").Value = 150000
$ws.Range("N9").Value = -150336
$ws.Range("H12").Value = 691.6667
$ws.Range("I12").Value = 475
$ws.Range("K12").Value = 475
$ws.Range("M12").Value = -305
$ws.Range("H15").Value = 25356.334
$ws.Range("J15").Value = 37500
$ws.Range("L15").Value = 37500
$ws.Range("N15").Value = -37840
$ws.Range("H76").Value = 5171
$ws.Range("I76").Value = 5171
$ws.Range("K76").Value = 5171
$ws.Range("M76").Value = -4856
$ws.Range("H79").Value = 5171
$ws.Range("I79").Value = 5171
$ws.Range("K79").Value = 5171
$ws.Range("M79").Value = -4079
$ws.Range("H107").Value = 2035.95
$ws.Range("I107").Value = 670.8889
$ws.Range("K107").Value = 670.8889
$ws.Range("M107").Value = 1249.1111
$ws.Range("H134").Value = 3464.6304
$ws.Range("I134").Value = 1718.36
$ws.Range("K134").Value = 5155.08
$ws.Range("M134").Value = -2620.08

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 404.53333
$ws.Range("I23").Value = 289.66666
$ws.Range("J23").Value = 481.1111
$ws.Range("K23").Value = 868.9999799999999
$ws.Range("L23").Value = 1443.3333
$ws.Range("M23").Value = -633.9999799999999
$ws.Range("N23").Value = -1913.3333
$ws.Range("H122").Value = 1489680.4
$ws.Range("J122").Value = 846.5333000000001
$ws.Range("L122").Value = 7618.7997
$ws.Range("N122").Value = -12518.7997
$ws.Range("H127").Value = 5050.5557
$ws.Range("J127").Value = 5050.5557
$ws.Range("L127").Value = 15151.6671
$ws.Range("N127").Value = -25071.6671
$ws.Range("H141").Value = 3986.1428
$ws.Range("I141").Value = 3986.1428
$ws.Range("K141").Value = 11958.4284
$ws.Range("M141").Value = -6778.428400000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H63").Value = 54897.75
$ws.Range("J63").Value = 54897.75
$ws.Range("L63").Value = 54897.75
$ws.Range("N63").Value = -56269.75
$ws.Range("H66").Value = 54897.75
$ws.Range("J66").Value = 54897.75
$ws.Range("L66").Value = 164693.25
$ws.Range("N66").Value = -171557.25
$ws.Range("H97").Value = 872.4167
$ws.Range("I97").Value = 1013.4167
$ws.Range("K97").Value = 1013.4167
$ws.Range("M97").Value = -517.4167

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5249.7
$ws.Range("I7").Value = 4473.875
$ws.Range("J7").Value = 5766.9165
$ws.Range("K7").Value = 4473.875
$ws.Range("L7").Value = 5766.9165
$ws.Range("M7").Value = -4361.875
$ws.Range("N7").Value = -5990.9165
$ws.Range("H16").Value = 315.1
$ws.Range("I16").Value = 257.375
$ws.Range("J16").Value = 546
$ws.Range("K16").Value = 257.375
$ws.Range("L16").Value = 546
$ws.Range("M16").Value = -87.375
$ws.Range("N16").Value = -886
$ws.Range("H46").Value = 5053411.5
$ws.Range("I46").Value = 2454.4546
$ws.Range("J46").Value = 10104368
$ws.Range("K46").Value = 2454.4546
$ws.Range("L46").Value = 10104368
$ws.Range("M46").Value = -2266.4546
$ws.Range("N46").Value = -10104744
$ws.Range("H50").Value = 50000
$ws.Range("J50").Value = 50000
$ws.Range("L50").Value = 50000
$ws.Range("N50").Value = -51274
$ws.Range("H93").Value = 929
$ws.Range("J93").Value = 807.5
$ws.Range("L93").Value = 807.5
$ws.Range("N93").Value = -3303.5
$ws.Range("H107").Value = 2966
$ws.Range("I107").Value = 2966
$ws.Range("K107").Value = 2966
$ws.Range("M107").Value = -1046
$ws.Range("H126").Value = 5249.7
$ws.Range("I126").Value = 4473.875
$ws.Range("J126").Value = 5766.9165
$ws.Range("K126").Value = 13421.625
$ws.Range("L126").Value = 17300.7495
$ws.Range("M126").Value = -10951.625
$ws.Range("N126").Value = -22240.7495
$ws.Range("H132").Value = 9811888
$ws.Range("I132").Value = 22730234
$ws.Range("K132").Value = 68190702
$ws.Range("M132").Value = -68188172

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 49296.332
$ws.Range("J64").Value = 49296.332
$ws.Range("L64").Value = 49296.332
$ws.Range("N64").Value = -49792.332
$ws.Range("H67").Value = 49296.332
$ws.Range("J67").Value = 49296.332
$ws.Range("L67").Value = 49296.332
$ws.Range("N67").Value = -51012.332
$ws.Range("H100").Value = 814.1429000000001
$ws.Range("I100").Value = 507.125
$ws.Range("J100").Value = 1223.5
$ws.Range("K100").Value = 1014.25
$ws.Range("L100").Value = 2447
$ws.Range("M100").Value = -473.25
$ws.Range("N100").Value = -3529
$ws.Range("H107").Value = 22223470
$ws.Range("I107").Value = 727
$ws.Range("K107").Value = 2181
$ws.Range("M107").Value = -261
$ws.Range("H132").Value = 3645.6428
$ws.Range("I132").Value = 3376.0857
$ws.Range("K132").Value = 10128.2571
$ws.Range("M132").Value = -7598.257100000001
